$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$s.Shapes | Get-Member | Out-String | Write-Output
